{"js": "// Replace the date and the 25 three-digit-by-one-digit multiplication\n// prompts with their updated values, mirroring the XML diff exactly.\n// Each \"old\" string is unique within the document, so a single\n// matchCase/wholeWord-less search+Replace per pair is safe.\n\nconst replacements = [\n  [\"2024-04-26 Friday\", \"2024-04-27 Saturday\"],\n  [\"513\u00d74=\", \"704\u00d77=\"],\n  [\"899\u00d77=\", \"867\u00d72=\"],\n  [\"394\u00d79=\", \"790\u00d73=\"],\n  [\"790\u00d72=\", \"678\u00d77=\"],\n  [\"707\u00d78=\", \"511\u00d72=\"],\n  [\"802\u00d75=\", \"525\u00d77=\"],\n  [\"528\u00d76=\", \"862\u00d78=\"],\n  [\"438\u00d72=\", \"838\u00d79=\"],\n  [\"978\u00d75=\", \"339\u00d79=\"],\n  [\"123\u00d74=\", \"819\u00d73=\"],\n  [\"702\u00d79=\", \"727\u00d76=\"],\n  [\"299\u00d74=\", \"505\u00d77=\"],\n  [\"200\u00d73=\", \"144\u00d72=\"],\n  [\"337\u00d73=\", \"655\u00d77=\"],\n  [\"651\u00d77=\", \"147\u00d79=\"],\n  [\"715\u00d78=\", \"182\u00d79=\"],\n  [\"510\u00d75=\", \"733\u00d74=\"],\n  [\"517\u00d73=\", \"413\u00d79=\"],\n  [\"518\u00d76=\", \"359\u00d79=\"],\n  [\"701\u00d77=\", \"894\u00d76=\"],\n  [\"539\u00d73=\", \"993\u00d76=\"],\n  [\"308\u00d77=\", \"329\u00d77=\"],\n  [\"916\u00d73=\", \"283\u00d74=\"],\n  [\"519\u00d75=\", \"786\u00d72=\"],\n  [\"816\u00d74=\", \"408\u00d78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date header and the 25 three-digit-by-one-digit\n# multiplication prompts to their new values (mirrors the XML diff 1:1).\n# Each \"old\" string occurs exactly once in the document, so a plain\n# Find/Replace (wdReplaceAll, match-case) per pair is exact and safe.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-04-26 Friday\", \"2024-04-27 Saturday\"),\n  @(\"513\u00d74=\", \"704\u00d77=\"),\n  @(\"899\u00d77=\", \"867\u00d72=\"),\n  @(\"394\u00d79=\", \"790\u00d73=\"),\n  @(\"790\u00d72=\", \"678\u00d77=\"),\n  @(\"707\u00d78=\", \"511\u00d72=\"),\n  @(\"802\u00d75=\", \"525\u00d77=\"),\n  @(\"528\u00d76=\", \"862\u00d78=\"),\n  @(\"438\u00d72=\", \"838\u00d79=\"),\n  @(\"978\u00d75=\", \"339\u00d79=\"),\n  @(\"123\u00d74=\", \"819\u00d73=\"),\n  @(\"702\u00d79=\", \"727\u00d76=\"),\n  @(\"299\u00d74=\", \"505\u00d77=\"),\n  @(\"200\u00d73=\", \"144\u00d72=\"),\n  @(\"337\u00d73=\", \"655\u00d77=\"),\n  @(\"651\u00d77=\", \"147\u00d79=\"),\n  @(\"715\u00d78=\", \"182\u00d79=\"),\n  @(\"510\u00d75=\", \"733\u00d74=\"),\n  @(\"517\u00d73=\", \"413\u00d79=\"),\n  @(\"518\u00d76=\", \"359\u00d79=\"),\n  @(\"701\u00d77=\", \"894\u00d76=\"),\n  @(\"539\u00d73=\", \"993\u00d76=\"),\n  @(\"308\u00d77=\", \"329\u00d77=\"),\n  @(\"916\u00d73=\", \"283\u00d74=\"),\n  @(\"519\u00d75=\", \"786\u00d72=\"),\n  @(\"816\u00d74=\", \"408\u00d78=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Replacement.ClearFormatting()\n  [void]$rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n\n"}
